$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the Hommes, JEL 2019/20 entry (row 12) as struck-through / superseded
$ws.Range("A12").Font.Strikethrough = $true
$ws.Range("B12").Font.Strikethrough = $true
$ws.Range("B12").WrapText = $true

# Add the new "Gobbi et al, 2019" record (SPF one-year ahead)
$ws.Range("A14").Value = "Gobbi et al, 2019"
$ws.Range("B14").Value = "the prob which agents assign to switching to liquidity trap regime is a metric for deanchoring. This prob, p, is determined via a logistic equation of the output gap. The model reconciles the empirical observations that 1) missing deflation + inflation, 2) de-anchoring of expectations and 3) steepening Phillips Curve (in terms of pi and x, in terms of pi and u, flattening). It also suggests that a Taylor rule may not be enough, even if aggressive, if shocks are big enough."
$ws.Range("B14").WrapText = $true
$ws.Rows(14).RowHeight = 45

$ws.Range("B23").Select()
